# Textbox response formatting fix
# Renames sheets and updates stimulus filename cells with new timestamp-based names.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (by positional index, matching workbook.xml order) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687425650241"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168745246005"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168745246434"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687453087118"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687453850935"

# --- Sheet 1 (GNG_TO) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651168742528181.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687425479298.csv"
$ws1.Range("B4").Value = "go_stims-1651168742548927.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687425640194.csv"

# --- Sheet 2 (NB_TO) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16511687429267974.csv"
$ws2.Range("B3").Value = "TB-16511687447908103.csv"
$ws2.Range("B4").Value = "ZB-match_4-16511687426641192.csv"
$ws2.Range("B5").Value = "TB-16511687452312796.csv"
$ws2.Range("B6").Value = "ZB-match_4-16511687428160548.csv"
$ws2.Range("B7").Value = "OB-1651168743677806.csv"
$ws2.Range("B8").Value = "ZB-match_3-16511687425818298.csv"
$ws2.Range("B9").Value = "OB-16511687431895618.csv"
$ws2.Range("B10").Value = "TB-1651168744494411.csv"

# --- Sheet 3 (RS_TO) --- no cell content changes, only the sheet name above

# --- Sheet 4 (TOL_TO) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687452612565.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687452494254.csv"
$ws4.Range("B4").Value = "MM_stims-16511687452926598.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687452612565.csv"
$ws4.Range("B6").Value = "MM_stims-16511687453087118.csv"
$ws4.Range("B7").Value = "ZM_stims-1651168745293659.csv"

# --- Sheet 5 (vSAT_TO) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16511687453695931.csv"
$ws5.Range("B3").Value = "SAT_stims-16511687453145788.csv"
$ws5.Range("B4").Value = "SAT_stims-16511687453401122.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651168745354115.csv"
